# Demo changes:
#  - Rename sheets in positions 4 and 5 ("xprodTestSheetName" / "prodTestSheetName")
#    so that the former becomes "prodTestSheetName" and the latter becomes
#    "xxprodTestSheetName". Rename sheet 5 first so the name "prodTestSheetName"
#    is freed up before sheet 4 claims it (Excel disallows duplicate sheet names).
#  - Make the (now) "prodTestSheetName" sheet (position 4) the active/selected tab,
#    instead of the sheet in position 5, which updates workbookView's activeTab
#    and moves tabSelected from the sheet 5 view to the sheet 4 view.

$wb = $excel.ActiveWorkbook

$sheet5 = $wb.Worksheets.Item(5)
$sheet5.Name = "xxprodTestSheetName"

$sheet4 = $wb.Worksheets.Item(4)
$sheet4.Name = "prodTestSheetName"

$sheet4.Activate()
